$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New blank "date style" cell at K2, cloned from the existing L2 format ---
$ws.Range("L2").Copy()
$ws.Range("K2").PasteSpecial(-4122)

# --- Header row: drop "username", shift the trailing headers left one column ---
$ws.Range("F1").Value = "password"
$ws.Range("G1").Value = "is_admin"
$ws.Range("H1").Value = "created_by"
$ws.Range("I1").Value = "updated_by"
$ws.Range("J1").Value = "remember_token"
$ws.Range("K1").Value = "created_at"
$ws.Range("L1").Value = "updated_at"
$ws.Range("M1").Clear()

# --- Addresses & phone numbers first (alamat / no_hp columns) ---
$ws.Range("C2").Value = "Jl.Prenjak Timur No 6, Sukun - Malang"
$ws.Range("C3").Value = "Jl.Madyopuro 6, Sawojajar"
$ws.Range("E2").Value = "'08871212"
$ws.Range("E3").Value = "'08871212"

# --- Row 2 (Yovie) ---
$ws.Range("B2").Value = "Yovie"
$ws.Range("F2").Value = 12345678
$ws.Range("G2").Value = 1
$ws.Range("H2").Clear()
$ws.Range("M2").Clear()

# --- Row 3 (Dewa) ---
$ws.Range("F3").Value = 12345678
$ws.Range("G3").Value = 0
$ws.Range("H3").Clear()

# --- Row 4 (Haykal) - new user ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Haykal"
$ws.Range("C4").Value = "Jl.Prenjak Timur No 6, Sukun - Malang"
$ws.Range("D4").Value = "haykal@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:haykal@gmail.com")
$ws.Range("D4").Style = "Hyperlink"
$ws.Range("E4").Value = "'08871212"
$ws.Range("F4").Value = 12345678
$ws.Range("G4").Value = 1

$ws.Range("I5").Select() | Out-Null
